$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11846
$ws.Range("C3").Value = 11846
$ws.Range("C4").Value = 11783
$ws.Range("C5").Value = 11409
$ws.Range("C6").Value = 10806
$ws.Range("C7").Value = 10806
$ws.Range("C8").Value = 10806
$ws.Range("C9").Value = 10806
$ws.Range("C10").Value = 10732
$ws.Range("C11").Value = 10689
$ws.Range("C12").Value = 9691
$ws.Range("C13").Value = 9691
$ws.Range("C14").Value = 9691
$ws.Range("C15").Value = 9691
$ws.Range("C16").Value = 9691
$ws.Range("C17").Value = 9691
$ws.Range("C18").Value = 9691
$ws.Range("C19").Value = 9691
$ws.Range("C20").Value = 9689
$ws.Range("C21").Value = 9364
$ws.Range("C22").Value = 9220
$ws.Range("C23").Value = 9220
$ws.Range("C24").Value = 8736
$ws.Range("C25").Value = 8736
$ws.Range("C26").Value = 8327
$ws.Range("C27").Value = 8327
$ws.Range("C28").Value = 8327
$ws.Range("C29").Value = 8269
$ws.Range("C30").Value = 8269
$ws.Range("C31").Value = 8262
$ws.Range("C32").Value = 8262
$ws.Range("C33").Value = 8262
$ws.Range("C34").Value = 8262
$ws.Range("C35").Value = 8262
$ws.Range("C36").Value = 8262
$ws.Range("C37").Value = 8262
$ws.Range("C38").Value = 7872
$ws.Range("C39").Value = 7872
$ws.Range("C40").Value = 7872
$ws.Range("C41").Value = 7872
$ws.Range("C42").Value = 7872
$ws.Range("C43").Value = 7872
$ws.Range("C44").Value = 7872
$ws.Range("C45").Value = 7872
$ws.Range("C46").Value = 7872
$ws.Range("C47").Value = 7870
$ws.Range("C48").Value = 7870
$ws.Range("C49").Value = 7870
$ws.Range("C50").Value = 7818
